{"js": "// Fix Race Collapse Variable: race_ethn_5cat previously erroneously omitted\n// NHAsian subjects (lumped into \"Other Race\"); restore the correct counts.\n//\n// Table layout (columns, 0-indexed):\n//   0 Variable | 1 Category | 2 N - TRUE | 3 Statistic (95% CI) - TRUE |\n//   4 N - FALSE | 5 Statistic (95% CI) - FALSE | 6 Test-Statistic (p-value)\n//\n// Rows (0-indexed) in the single table of this document:\n//   7 race_ethn_5cat / Non-Hispanic Asian\n//   8 race_ethn_5cat / Other Race (including multiracial)\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load(\"rowCount, values\");\nawait context.sync();\n\n// Locate the two target rows by their Variable/Category labels instead of\n// hard-coded indices, so the script is resilient if row order ever shifts.\nconst values = table.values;\n\nfunction findRow(variable, category) {\n  for (let r = 0; r < values.length; r++) {\n    const row = values[r];\n    if (row[0] === variable && row[1] === category) {\n      return r;\n    }\n  }\n  throw new Error(`Row not found for ${variable} / ${category}`);\n}\n\nconst asianRow = findRow(\"race_ethn_5cat\", \"Non-Hispanic Asian\");\nconst otherRow = findRow(\"race_ethn_5cat\", \"Other Race (including multiracial)\");\n\nconst edits = [\n  // row, col, newText\n  [asianRow, 2, \"87\"],\n  [asianRow, 3, \"11.3 (9.3 - 13.8)\"],\n  [asianRow, 4, \"37\"],\n  [asianRow, 5, \"9.2 (6.7 - 12.5)\"],\n  [otherRow, 2, \"36\"],\n  [otherRow, 3, \"4.7 (3.4 - 6.4)\"],\n  [otherRow, 4, \"12\"],\n  [otherRow, 5, \"3.0 (1.7 - 5.2)\"],\n];\n\nfor (const [row, col, text] of edits) {\n  const cell = table.getCell(row, col);\n  cell.value = text;\n}\n\nawait context.sync();\n", "ps1": "# Fix Race Collapse Variable: race_ethn_5cat previously erroneously omitted\n# NHAsian subjects (lumped into \"Other Race\"); restore the correct counts.\n#\n# Table layout (columns, 1-indexed, as used by the Word COM object model):\n#   1 Variable | 2 Category | 3 N - TRUE | 4 Statistic (95% CI) - TRUE |\n#   5 N - FALSE | 6 Statistic (95% CI) - FALSE | 7 Test-Statistic (p-value)\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\nfunction Get-CellText($table, $r, $c) {\n    return $table.Cell($r, $c).Range.Text.TrimEnd([char]13, [char]7)\n}\n\nfunction Set-CellText($table, $r, $c, $text) {\n    $table.Cell($r, $c).Range.Text = $text\n}\n\nfunction Find-RowByVariableCategory($table, $variable, $category) {\n    for ($r = 1; $r -le $table.Rows.Count; $r++) {\n        $v = Get-CellText $table $r 1\n        $cat = Get-CellText $table $r 2\n        if ($v -eq $variable -and $cat -eq $category) {\n            return $r\n        }\n    }\n    return $null\n}\n\n$asianRow = Find-RowByVariableCategory $t \"race_ethn_5cat\" \"Non-Hispanic Asian\"\n$otherRow = Find-RowByVariableCategory $t \"race_ethn_5cat\" \"Other Race (including multiracial)\"\n\nSet-CellText $t $asianRow 3 \"87\"\nSet-CellText $t $asianRow 4 \"11.3 (9.3 - 13.8)\"\nSet-CellText $t $asianRow 5 \"37\"\nSet-CellText $t $asianRow 6 \"9.2 (6.7 - 12.5)\"\n\nSet-CellText $t $otherRow 3 \"36\"\nSet-CellText $t $otherRow 4 \"4.7 (3.4 - 6.4)\"\nSet-CellText $t $otherRow 5 \"12\"\nSet-CellText $t $otherRow 6 \"3.0 (1.7 - 5.2)\"\n"}
